{"js": "// Replace each two-digit multiplication prompt with its new value.\n// Mapping derived from the canonical OOXML diff (old => new), applied\n// in document order; each old string is unique within the document so\n// a plain text search/replace is unambiguous.\nconst replacements = [\n  [\"20\u00d722=\", \"93\u00d750=\"],\n  [\"67\u00d791=\", \"12\u00d729=\"],\n  [\"20\u00d796=\", \"58\u00d790=\"],\n  [\"29\u00d725=\", \"13\u00d741=\"],\n  [\"98\u00d714=\", \"54\u00d713=\"],\n  [\"20\u00d785=\", \"90\u00d724=\"],\n  [\"67\u00d724=\", \"74\u00d718=\"],\n  [\"15\u00d795=\", \"52\u00d759=\"],\n  [\"66\u00d711=\", \"94\u00d740=\"],\n  [\"47\u00d716=\", \"15\u00d741=\"],\n  [\"28\u00d730=\", \"33\u00d784=\"],\n  [\"21\u00d712=\", \"82\u00d718=\"],\n  [\"69\u00d739=\", \"40\u00d760=\"],\n  [\"79\u00d757=\", \"36\u00d754=\"],\n  [\"65\u00d768=\", \"89\u00d784=\"],\n  [\"38\u00d763=\", \"19\u00d753=\"],\n  [\"36\u00d713=\", \"11\u00d718=\"],\n  [\"33\u00d759=\", \"63\u00d766=\"],\n  [\"53\u00d767=\", \"82\u00d747=\"],\n  [\"61\u00d798=\", \"89\u00d751=\"],\n  [\"91\u00d736=\", \"74\u00d718=\"],\n  [\"56\u00d717=\", \"85\u00d746=\"],\n  [\"25\u00d743=\", \"28\u00d783=\"],\n  [\"32\u00d751=\", \"17\u00d767=\"],\n  [\"14\u00d738=\", \"18\u00d778=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication prompt with its new value.\n# Mapping derived from the canonical OOXML diff (old => new); each old\n# string is unique in the document, so Find/Replace (replace-all) is\n# unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old = \"20\u00d722=\"; new = \"93\u00d750=\"},\n    @{old = \"67\u00d791=\"; new = \"12\u00d729=\"},\n    @{old = \"20\u00d796=\"; new = \"58\u00d790=\"},\n    @{old = \"29\u00d725=\"; new = \"13\u00d741=\"},\n    @{old = \"98\u00d714=\"; new = \"54\u00d713=\"},\n    @{old = \"20\u00d785=\"; new = \"90\u00d724=\"},\n    @{old = \"67\u00d724=\"; new = \"74\u00d718=\"},\n    @{old = \"15\u00d795=\"; new = \"52\u00d759=\"},\n    @{old = \"66\u00d711=\"; new = \"94\u00d740=\"},\n    @{old = \"47\u00d716=\"; new = \"15\u00d741=\"},\n    @{old = \"28\u00d730=\"; new = \"33\u00d784=\"},\n    @{old = \"21\u00d712=\"; new = \"82\u00d718=\"},\n    @{old = \"69\u00d739=\"; new = \"40\u00d760=\"},\n    @{old = \"79\u00d757=\"; new = \"36\u00d754=\"},\n    @{old = \"65\u00d768=\"; new = \"89\u00d784=\"},\n    @{old = \"38\u00d763=\"; new = \"19\u00d753=\"},\n    @{old = \"36\u00d713=\"; new = \"11\u00d718=\"},\n    @{old = \"33\u00d759=\"; new = \"63\u00d766=\"},\n    @{old = \"53\u00d767=\"; new = \"82\u00d747=\"},\n    @{old = \"61\u00d798=\"; new = \"89\u00d751=\"},\n    @{old = \"91\u00d736=\"; new = \"74\u00d718=\"},\n    @{old = \"56\u00d717=\"; new = \"85\u00d746=\"},\n    @{old = \"25\u00d743=\"; new = \"28\u00d783=\"},\n    @{old = \"32\u00d751=\"; new = \"17\u00d767=\"},\n    @{old = \"14\u00d738=\"; new = \"18\u00d778=\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.old\n    $find.Replacement.Text = $r.new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n}\n"}
